$wb = $excel.ActiveWorkbook

# Rename sheet tabs
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778309508793"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778325629091"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778325638773"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778326278872"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650477832690912"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778309198797.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778309339118.csv"
$ws1.Range("B4").Value = "go_stims-16504778309348788.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778309499109.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16504778325398753.csv"
$ws2.Range("B3").Value = "ZB-match_9-16504778313088818.csv"
$ws2.Range("B4").Value = "TB-1650477832430879.csv"
$ws2.Range("B5").Value = "OB-16504778313888755.csv"
$ws2.Range("B6").Value = "OB-16504778315388782.csv"
$ws2.Range("B7").Value = "ZB-match_1-16504778309768796.csv"
$ws2.Range("B8").Value = "OB-1650477832011881.csv"
$ws2.Range("B9").Value = "ZB-match_4-1650477831252879.csv"
$ws2.Range("B10").Value = "TB-16504778324558785.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650477832594912.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778325698786.csv"
$ws4.Range("B4").Value = "MM_stims-16504778326109114.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477832594912.csv"
$ws4.Range("B6").Value = "MM_stims-16504778326268766.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778326109114.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778326589108.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778326749113.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778326429107.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778326308813.csv"
